# Adds the team's season record (Wins / Losses / Ties) as three new
# trailing columns (AD, AE, AF) to the roster sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Season record for this team/year.
$wins   = 68
$losses = 94
$ties   = 0

# New header cells, styled like the rest of row 1 (bold + border, via the
# same "header" format already used by the last existing header cell).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill the record for every player row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 69
}

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}
